# Update BOM with extra caps:
#  - Quantity for C13, C14, C15, C36 (row 31) increases from 4 to 6
#  - Extended Price column (J2:J49) is cleared out

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump quantity for the capacitor row (C13, C14, C15, C36) from 4 to 6
$ws.Range("H31").Value = 6

# Clear the stale "Extended Price" values for every BOM line
$ws.Range("J2:J49").ClearContents()

# Leave the selection where the editor left it after the delete
$ws.Range("J2:J49").Select() | Out-Null
